$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 64).Value = 0.93112304316550421
$ws.Cells.Item(3, 52).Value = 0.97781445578991866
$ws.Cells.Item(4, 2).Value = 0.88896679932301148
$ws.Cells.Item(4, 16).Value = 0.98722566806428325
$ws.Cells.Item(4, 39).Value = 0.78014933551683208
$ws.Cells.Item(6, 5).Value = 0.76305515147283964
$ws.Cells.Item(6, 41).Value = 0.72677986819599649
$ws.Cells.Item(6, 66).Value = 0.74696672407579334
$ws.Cells.Item(7, 5).Value = 0.93396921510692132
$ws.Cells.Item(7, 40).Value = 0.84338332880934264
$ws.Cells.Item(8, 57).Value = 0.80855040175186299
$ws.Cells.Item(8, 60).Value = 0.97331476474474521
$ws.Cells.Item(8, 66).Value = 0.81197210880449244
$ws.Cells.Item(9, 11).Value = 0.88171512286065767
$ws.Cells.Item(10, 9).Value = 0.62552645242119298
$ws.Cells.Item(10, 30).Value = 0.86572850922808453
$ws.Cells.Item(10, 44).Value = 0.79659284106791117
$ws.Cells.Item(10, 65).Value = 0.89001800549693189
$ws.Cells.Item(11, 16).Value = 0.7371151500819173
$ws.Cells.Item(11, 29).Value = 0.98913332097955475
$ws.Cells.Item(11, 58).Value = 0.98259003407373757
$ws.Cells.Item(12, 24).Value = 0.56705440106061533
$ws.Cells.Item(12, 38).Value = 0.78568168515511583
$ws.Cells.Item(13, 15).Value = 0.91322682173115499
$ws.Cells.Item(14, 13).Value = 0.97492063412813057
$ws.Cells.Item(14, 25).Value = 0.86315004877144452
$ws.Cells.Item(14, 34).Value = 0.93228364982335965
$ws.Cells.Item(14, 43).Value = 0.92783887207918347
$ws.Cells.Item(14, 52).Value = 0.88926217756996984
$ws.Cells.Item(15, 6).Value = 0.91641898090903662
$ws.Cells.Item(16, 33).Value = 0.95696829979695874
$ws.Cells.Item(17, 24).Value = 0.85583014515503208
$ws.Cells.Item(17, 32).Value = 0.90766890913946607
$ws.Cells.Item(17, 45).Value = 0.84644450409673022
$ws.Cells.Item(17, 49).Value = 0.90172280610893873
$ws.Cells.Item(17, 62).Value = 0.78927662136523735
$ws.Cells.Item(18, 11).Value = 0.99941717979435918
$ws.Cells.Item(18, 36).Value = 0.85578485128015247
$ws.Cells.Item(18, 51).Value = 0.61768370740677603
$ws.Cells.Item(19, 20).Value = 0.73190665571107738
$ws.Cells.Item(19, 32).Value = 0.99179179257844974
$ws.Cells.Item(19, 35).Value = 0.72613828633020683
$ws.Cells.Item(19, 39).Value = 0.7365111348778508
$ws.Cells.Item(20, 29).Value = 0.90282900868435534
$ws.Cells.Item(20, 48).Value = 0.9260305453496217
$ws.Cells.Item(20, 57).Value = 0.77982282446548634
$ws.Cells.Item(20, 68).Value = 0.72679906446907938
$ws.Cells.Item(21, 1).Value = 0.94202211866494268
$ws.Cells.Item(22, 21).Value = 0.89886719480029376
$ws.Cells.Item(22, 39).Value = 0.98657091074248571
$ws.Cells.Item(23, 49).Value = 0.91025629807658781
$ws.Cells.Item(25, 9).Value = 0.99424867758771174
$ws.Cells.Item(26, 15).Value = 0.99944778795196421
$ws.Cells.Item(26, 24).Value = 0.84429727529768228
$ws.Cells.Item(26, 62).Value = 0.71017466742445401
$ws.Cells.Item(27, 6).Value = 0.76507527477749249
$ws.Cells.Item(27, 7).Value = 0.79020654709520821
$ws.Cells.Item(27, 26).Value = 0.59191804423401417
$ws.Cells.Item(27, 28).Value = 0.7476389344312514
$ws.Cells.Item(28, 3).Value = 0.95571332099419637
$ws.Cells.Item(29, 59).Value = 0.84919499031400414
$ws.Cells.Item(30, 5).Value = 0.64816968618409709
$ws.Cells.Item(30, 28).Value = 0.93017937252603744
$ws.Cells.Item(30, 51).Value = 0.8203217031370893
$ws.Cells.Item(30, 65).Value = 0.91480056116024044
$ws.Cells.Item(31, 2).Value = 0.73207207710584266
$ws.Cells.Item(31, 22).Value = 0.73195424345242666
$ws.Cells.Item(33, 35).Value = 0.89682770003949508
$ws.Cells.Item(35, 31).Value = 0.92772885824447093
$ws.Cells.Item(35, 65).Value = 0.8823826863645533
$ws.Cells.Item(36, 12).Value = 0.97542171506956055
$ws.Cells.Item(36, 37).Value = 0.70905680460624032
$ws.Cells.Item(37, 28).Value = 0.9218346883534444
$ws.Cells.Item(38, 35).Value = 0.67152600145419572
$ws.Cells.Item(39, 2).Value = 0.59012407062111694
$ws.Cells.Item(39, 67).Value = 0.9272480179187923
$ws.Cells.Item(40, 5).Value = 0.84347879359717193
$ws.Cells.Item(40, 68).Value = 0.83375107870872056
$ws.Cells.Item(41, 5).Value = 0.66068185392881307
$ws.Cells.Item(42, 16).Value = 0.89483868384772602
$ws.Cells.Item(42, 32).Value = 0.87092902317907561
$ws.Cells.Item(43, 4).Value = 0.75276587196619937
$ws.Cells.Item(43, 45).Value = 0.93839373015009253
$ws.Cells.Item(44, 30).Value = 0.89569848838221888
$ws.Cells.Item(45, 33).Value = 0.74323999639524274
$ws.Cells.Item(45, 37).Value = 0.94546383694095804
$ws.Cells.Item(46, 45).Value = 0.78806849554266667
$ws.Cells.Item(46, 48).Value = 0.77341505529474541
$ws.Cells.Item(46, 63).Value = 0.7122103602533647
$ws.Cells.Item(47, 18).Value = 0.73807095101436648
$ws.Cells.Item(47, 36).Value = 0.89824138693794731
$ws.Cells.Item(48, 25).Value = 0.83562564717629884
$ws.Cells.Item(48, 44).Value = 0.62979529671568768
$ws.Cells.Item(48, 50).Value = 0.80006285118665477
$ws.Cells.Item(49, 1).Value = 0.95955220516618689
$ws.Cells.Item(50, 9).Value = 0.78322518528342822
$ws.Cells.Item(50, 49).Value = 0.68148510321386135
$ws.Cells.Item(50, 54).Value = 0.78677353709325404
$ws.Cells.Item(51, 16).Value = 0.58395185102930047
$ws.Cells.Item(52, 28).Value = 0.74922979494080266
$ws.Cells.Item(52, 34).Value = 0.82984231444441936
$ws.Cells.Item(52, 49).Value = 0.73256076574355145
$ws.Cells.Item(53, 12).Value = 0.87844048559718679
$ws.Cells.Item(53, 23).Value = 0.79989750167024409
$ws.Cells.Item(53, 39).Value = 0.54513837268619048
$ws.Cells.Item(53, 54).Value = 0.94249287457598885
$ws.Cells.Item(54, 49).Value = 0.92611901543799324
$ws.Cells.Item(54, 51).Value = 0.87393206463670725
$ws.Cells.Item(55, 23).Value = 0.9268200333648684
$ws.Cells.Item(55, 48).Value = 0.98668796334649644
$ws.Cells.Item(56, 7).Value = 0.96461775084589385
$ws.Cells.Item(56, 54).Value = 0.70613320509323008
$ws.Cells.Item(57, 30).Value = 0.94059093735251609
$ws.Cells.Item(57, 64).Value = 0.95182001172073405
$ws.Cells.Item(58, 13).Value = 0.86283886259741172
$ws.Cells.Item(58, 53).Value = 0.72403072719563544
$ws.Cells.Item(58, 56).Value = 0.98314804398238209
$ws.Cells.Item(59, 54).Value = 0.72361057969268883
$ws.Cells.Item(60, 52).Value = 0.95492099384731555
$ws.Cells.Item(61, 20).Value = 0.78482418141589028
$ws.Cells.Item(61, 21).Value = 0.83258225316644729
$ws.Cells.Item(61, 63).Value = 0.93147600533948571
$ws.Cells.Item(62, 13).Value = 0.72957523222108578
$ws.Cells.Item(62, 32).Value = 0.89995736917261271
$ws.Cells.Item(62, 65).Value = 0.77300045597648293
$ws.Cells.Item(63, 31).Value = 0.65658089586253732
$ws.Cells.Item(63, 38).Value = 0.96768652579007908
$ws.Cells.Item(63, 64).Value = 0.69451004686053419
$ws.Cells.Item(64, 61).Value = 0.67165699691526592
$ws.Cells.Item(65, 34).Value = 0.9355261263171869
$ws.Cells.Item(65, 39).Value = 0.96918026281993885
$ws.Cells.Item(66, 7).Value = 0.8167702840987876
$ws.Cells.Item(66, 11).Value = 0.90061459239350561
$ws.Cells.Item(66, 36).Value = 0.77548639141734443
$ws.Cells.Item(67, 19).Value = 0.6445562889109262
$ws.Cells.Item(67, 60).Value = 0.90312665591189889
